# Update the "取得日時" (acquisition timestamp) column on the "ランサーズ" sheet
# for rows 2-8 to reflect the new run time: 2025-12-20 18:31:13

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-12-20 18:31:13"

for ($row = 2; $row -le 8; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
